# GSC export refresh: the rolling date window advanced by one day.
# The oldest date row (2025-11-13) drops off the front of the "Chart"
# sheet's data table and every subsequent row shifts up to take its place
# (Excel's row-delete semantics also renumber the now-unused shared
# string and keep the "Critical issues"/"Non-critical issues" header
# rows' string references consistent automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the oldest date (2025-11-13, directly under the header row).
# Deleting it shifts all the remaining date rows up by one.
$ws.Rows.Item(2).Delete()
